# Generate Report for Handoff
# Adds two new "Ready for handoff" entries (87493161-... and eea7592c-...)
# to the Overview sheet and to each language sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$uuid1 = "87493161-28ba-4f50-9975-e968af648f58"
$uuid2 = "eea7592c-613b-4bca-be7d-0cbe7649faeb"

$hash1 = "27432036023a1cf0aec99ef47635d5dd41f17df2"
$hash2 = "935cd65346599f4555ea17e579e6d3c451a85d7d"

$status = "Ready for handoff"
$dtFormat = "yyyy-mm-dd HH:mm:ss"
$neverHandedBack = "0001-01-01 00:00:00"

$overviewDate = "2016-03-23 16:43:16"
$zhcnHandoffDatetime = "2016-03-23 16:43:12"
$dedeHandoffDatetime = "2016-03-23 16:43:16"

# ============================================================
# Sheet "Overview"
# ============================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "$uuid1.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid1.md", [Type]::Missing, [Type]::Missing, "$uuid1.md")
$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status
$wsOverview.Range("D4").Value = $overviewDate
$wsOverview.Range("D4").NumberFormat = $dtFormat

$wsOverview.Range("A5").Value = "$uuid2.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid2.md", [Type]::Missing, [Type]::Missing, "$uuid2.md")
$wsOverview.Range("B5").Value = $status
$wsOverview.Range("C5").Value = $status
$wsOverview.Range("D5").Value = $overviewDate
$wsOverview.Range("D5").NumberFormat = $dtFormat

# ============================================================
# Sheet "zh-cn"
# ============================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlf1 = "$uuid1.$hash1.zh-cn.xlf"
$zhXlf2 = "$uuid2.$hash2.zh-cn.xlf"

$wsZhCn.Range("A4").Value = "$uuid1.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid1.md", [Type]::Missing, [Type]::Missing, "$uuid1.md")
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = $status
$wsZhCn.Range("D4").Value = $zhXlf1
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf1", [Type]::Missing, [Type]::Missing, $zhXlf1)
$wsZhCn.Range("E4").Value = $zhcnHandoffDatetime
$wsZhCn.Range("E4").NumberFormat = $dtFormat
$wsZhCn.Range("H4").Value = $neverHandedBack
$wsZhCn.Range("H4").NumberFormat = $dtFormat
$wsZhCn.Range("J4").Value = "Include"

$wsZhCn.Range("A5").Value = "$uuid2.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid2.md", [Type]::Missing, [Type]::Missing, "$uuid2.md")
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = $status
$wsZhCn.Range("D5").Value = $zhXlf2
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf2", [Type]::Missing, [Type]::Missing, $zhXlf2)
$wsZhCn.Range("E5").Value = $zhcnHandoffDatetime
$wsZhCn.Range("E5").NumberFormat = $dtFormat
$wsZhCn.Range("H5").Value = $neverHandedBack
$wsZhCn.Range("H5").NumberFormat = $dtFormat
$wsZhCn.Range("J5").Value = "Include"

# ============================================================
# Sheet "de-de"
# ============================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlf1 = "$uuid1.$hash1.de-de.xlf"
$deXlf2 = "$uuid2.$hash2.de-de.xlf"

$wsDeDe.Range("A4").Value = "$uuid1.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid1.md", [Type]::Missing, [Type]::Missing, "$uuid1.md")
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = $status
$wsDeDe.Range("D4").Value = $deXlf1
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf1", [Type]::Missing, [Type]::Missing, $deXlf1)
$wsDeDe.Range("E4").Value = $dedeHandoffDatetime
$wsDeDe.Range("E4").NumberFormat = $dtFormat
$wsDeDe.Range("H4").Value = $neverHandedBack
$wsDeDe.Range("H4").NumberFormat = $dtFormat
$wsDeDe.Range("J4").Value = "Include"

$wsDeDe.Range("A5").Value = "$uuid2.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$uuid2.md", [Type]::Missing, [Type]::Missing, "$uuid2.md")
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = $status
$wsDeDe.Range("D5").Value = $deXlf2
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf2", [Type]::Missing, [Type]::Missing, $deXlf2)
$wsDeDe.Range("E5").Value = $dedeHandoffDatetime
$wsDeDe.Range("E5").NumberFormat = $dtFormat
$wsDeDe.Range("H5").Value = $neverHandedBack
$wsDeDe.Range("H5").NumberFormat = $dtFormat
$wsDeDe.Range("J5").Value = "Include"

Write-Output "Added handoff rows for $uuid1 and $uuid2 to Overview, zh-cn, and de-de sheets."
